$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")
$ws.Activate()

# Update the repaymentstrategy value (row 17, column B) from "RBI (India)"
# to the new scenario value "Overdue/Due Fee/Int,Principal"
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Move the active selection to the edited cell, matching the authored change
$ws.Range("B17").Select()
